# Update cryptocurrency price/volume snapshot per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.631.67'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").Value = '''1.824.22'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("D4").Value = '''1.008'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''1.008'
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = '''307.59'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("E7").Value = '  +2.33%  '
$ws.Range("D8").Value = '''0.3605'
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").Value = '''0.07136'
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").Value = '''0.9005'
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").Value = '''0.07765'
$ws.Range("E11").Value = '  -0.68%  '
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").Value = '''1.793.07'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = '''5.261'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '''87.71'
$ws.Range("E16").Value = '  +3.39%  '
$ws.Range("D17").Value = '''1.010'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '''0.000008570'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '''26.679.93'
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '''14.15'
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '''1.920'
$ws.Range("E24").Value = '  -2.37%  '
$ws.Range("D25").Value = '''152.29'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = '''17.90'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''1.969'
$ws.Range("E27").Value = '  -3.44%  '
$ws.Range("D28").Value = '''113.79'
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("D29").Value = '''4.821'
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("D30").Value = '''0.08807'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("D31").Value = '''3.140'
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("D32").Value = '''0.7313'
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.434'
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '''1.137'
$ws.Range("E34").Value = '  +2.42%  '
$ws.Range("D35").Value = '''2.699'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D37").Value = '''0.01923'
$ws.Range("E37").Value = '  -0.59%  '
$ws.Range("D38").Value = '''2.925'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").Value = '''0.05112'
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").Value = '''6.895'
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("D41").Value = '''0.5055'
$ws.Range("E41").Value = '  -1.54%  '
$ws.Range("E42").Value = '  -1.57%  '
$ws.Range("D43").Value = '''7.993'
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '''0.4657'
$ws.Range("E45").Value = '  +0.12%  '
$ws.Range("D46").Value = '''9.963'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("D47").Value = '''98.35'
$ws.Range("E47").Value = '  -2.38%  '
$ws.Range("D48").Value = '''1.558'
$ws.Range("E48").Value = '  -1.64%  '
$ws.Range("D49").Value = '''0.05986'
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").Value = '''35.84'
$ws.Range("E51").Value = '  -0.93%  '
